# Generate Report for Handback
# Update the handback status report timestamps/status produced by a new
# report run (new Xliff generate / handoff / handback datetimes, and the
# zh-cn Priority flipping from human-translation "ht" to machine-translation
# "mt").

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" column (G) for rows 2 and 4
$wsOverview.Range("G2").Value = "2016-08-26 22:16:35"
$wsOverview.Range("G4").Value = "2016-08-26 22:16:35"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# "Priority" column (E) for rows 2 and 4: human translation -> machine translation
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# "Correspond Handoff Datetime" column (H) for rows 2 and 4
$wsZhCn.Range("H2").Value = "2016-08-26 22:16:29"
$wsZhCn.Range("H4").Value = "2016-08-26 22:16:29"
# "Correspond Handback DateTime" column (K) for rows 2 and 4
$wsZhCn.Range("K2").Value = "2016-08-26 22:16:56"
$wsZhCn.Range("K4").Value = "2016-08-26 22:16:56"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# "Correspond Handoff Datetime" column (H) for rows 2 and 4
$wsDeDe.Range("H2").Value = "2016-08-26 22:16:35"
$wsDeDe.Range("H4").Value = "2016-08-26 22:16:35"
# "Correspond Handback DateTime" column (K) for rows 2 and 4
$wsDeDe.Range("K2").Value = "2016-08-26 22:17:08"
$wsDeDe.Range("K4").Value = "2016-08-26 22:17:08"
